$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.031.19'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.548.92'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.56%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '538.74'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.12'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.572'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.572.89'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.48%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.54%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.48'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.363'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.000.15'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '24.07'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '59.992.62'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.77%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000143'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.535.16'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.63%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.80%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '326.98'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.98'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.33'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.25%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.53%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.995'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.03'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.05'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0796'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.03%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.19'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -5.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '165.12'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.20%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.02%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.74'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.09%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.02'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '301.65'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.64%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.58'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -6.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.835'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +5.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.73'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.00%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.610'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.59%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '127.35'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.00%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0939'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.95'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.07%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0228'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.48%  '
